$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.233.17"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "2.432.31"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "489.95"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.43"
$ws.Range("E6").Value = "  +3.44%  "
$ws.Range("E7").Value = "  +20.34%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "2.453.15"
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.66"
$ws.Range("E11").Value = "  -2.20%  "
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "2.857.02"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").Value = "57.240.02"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.88"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "2.449.08"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.80"
$ws.Range("E19").Value = "  +5.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "329.60"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.00"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.95"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.32"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.412"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "2.543.75"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").Value = "  -4.02%  "
$ws.Range("D30").Value = "0.0₃0793"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.77"
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.96"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.33"
$ws.Range("E35").Value = "  +2.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.16"
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.73"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.858"
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("E39").Value = "  +11.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.27"
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.39"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.54"
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.600"
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("E45").Value = "  -3.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "267.80"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.23"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("E49").Value = "  -4.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.58"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.75"
$ws.Range("E51").Value = "  +14.96%  "
